$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new export timestamp
$ws.Name = "牌局記錄0810_013525"

# New play-history rows to append (columns A:L), matching the source data
$data = @(
    @(3,'J,Q,3',1,'Q,2,9',2,'閒','紅','藍','藍','藍','紅','紅'),
    @(2,'9,5,8',4,'5,9',2,'莊','藍','藍','紅','紅','紅','藍'),
    @(0,'K,10,Q',6,'1,Q,5',6,'莊','紅','藍','紅','藍','紅','藍'),
    @(4,'2,2',9,'9,Q',5,'莊','藍','紅','紅','紅','藍','藍'),
    @(9,'3,6',4,'10,4',5,'閒','紅','紅','藍','藍','藍','紅'),
    @(3,'10,4,9',8,'K,J,8',5,'莊','藍','藍','紅','紅','紅','藍'),
    @(4,'4,K',8,'8,Q',4,'莊','紅','紅','紅','藍','藍','藍'),
    @(7,'7,K',9,'J,9',2,'莊','藍','藍','藍','紅','紅','紅'),
    @(8,'3,5',0,'K,10',8,'閒','紅','紅','紅','藍','藍','藍'),
    @(8,'4,6,8',4,'Q,10,4',4,'閒','藍','紅','紅','紅','藍','藍'),
    @(0,'J,Q,K',1,'K,3,8',1,'莊','紅','藍','藍','藍','紅','紅'),
    @(1,'8,7,6',3,'1,2,Q',2,'莊','紅','紅','藍','藍','藍','紅'),
    @(5,'9,6',9,'9,10',4,'莊','藍','藍','紅','紅','紅','藍'),
    @(7,'1,6',5,'7,6,2',2,'閒','紅','紅','紅','藍','藍','藍'),
    @(1,'5,6,K',2,'J,2,Q',1,'莊','藍','藍','藍','紅','紅','紅'),
    @(0,'3,K,7',9,'5,5,9',9,'莊','藍','紅','紅','紅','藍','藍'),
    @(2,'7,5',9,'Q,9',7,'莊','紅','紅','藍','藍','藍','紅'),
    @(5,'J,5',8,'6,2',3,'莊','紅','藍','紅','藍','紅','藍')
)

$startRow = 900
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
